# mm2-diagrams.pptx -- "work on test tool doc"
#
# Slide 20 ("Test Environment: Bidirectional") renames the Elasticsearch
# index label "assets" to "accounts" in four textboxes that annotate the
# MirrorMaker 2.0 replication diagram. Each textbox auto-fits its width to
# the text (wrap="none" + spAutoFit), so PowerPoint also widens the shape
# to keep the label from clipping. The widths below are the exact point
# values that round-trip to the target EMU widths (971741 / 694421 /
# 694421 / 893193) through this host's single-precision Width property.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(20)

# "esm.assets" -> "esm.accounts" (left textbox, near the esm cluster)
$sh = $s.Shapes.Item(21)
$sh.TextFrame.TextRange.Text = "esm.accounts"
$sh.Width = 76.51504516601562

# "assets" -> "accounts" (textbox on the esm-to-esi connector)
$sh = $s.Shapes.Item(22)
$sh.TextFrame.TextRange.Text = "accounts"
$sh.Width = 54.6788215637207

# "assets" -> "accounts" (textbox on the esi-to-esm connector)
$sh = $s.Shapes.Item(28)
$sh.TextFrame.TextRange.Text = "accounts"
$sh.Width = 54.6788215637207

# "esi.assets" -> "esi.accounts" (right textbox, near the esi cluster)
$sh = $s.Shapes.Item(29)
$sh.TextFrame.TextRange.Text = "esi.accounts"
$sh.Width = 70.33016204833984
